$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at 482, shifting existing rows 482:554 down to 483:555.
$ws.Rows.Item(482).Insert()

# Populate the newly-inserted row 482 with the new price entry.
$ws.Cells.Item(482, 1).Value = 6
$ws.Cells.Item(482, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(482, 3).Value = "Metropolitana"
$ws.Cells.Item(482, 4).Value = 44776
$ws.Cells.Item(482, 5).Value = 13
$ws.Cells.Item(482, 6).Value = 100112044
$ws.Cells.Item(482, 7).Value = "Perejil"
$ws.Cells.Item(482, 8).Value = "Sin especificar"
$ws.Cells.Item(482, 9).Value = "Primera"
$ws.Cells.Item(482, 10).Value = 150
$ws.Cells.Item(482, 11).Value = 19000
$ws.Cells.Item(482, 12).Value = 20000
$ws.Cells.Item(482, 13).Value = 19600
$ws.Cells.Item(482, 14).Value = "`$/docena de atados"
$ws.Cells.Item(482, 15).Value = "Región Metropolitana"
$ws.Cells.Item(482, 16).Value = 6533
$ws.Cells.Item(482, 17).Value = 3
$ws.Cells.Item(482, 18).Value = "Hortaliza"
